$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# Helper: find the single paragraph whose trimmed visible text matches $needle
# exactly (paragraph.Range.Text does not surface text that lives inside
# w:smartTag elements in this runtime, so the comparisons below only use the
# text that sits outside of smart tags).
function Get-ParaByExactText($doc, [string]$needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text.Trim() -eq $needle) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1) "Date: <Sept 1, 2008> - <Sept 3, 2008>" paragraph (near top of doc):
#    - reorder the w:attr children of both smartTagPr elements from
#      Month/Day/Year to Year/Day/Month.
# ---------------------------------------------------------------------------
$pDate = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.Trim().StartsWith("Date:")) {
        $pDate = $p
        break
    }
}
$xmlDate = @'
<w:p w:rsidR="00CF03FF" w:rsidRPr="00D37CDD" w:rsidRDefault="00CF03FF" w:rsidP="00FF5F4B"><w:pPr><w:ind w:left="284"/><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr></w:pPr><w:r w:rsidRPr="00D37CDD"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t xml:space="preserve">Date: </w:t></w:r><w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date"><w:smartTagPr><w:attr w:name="Year" w:val="2008"/><w:attr w:name="Day" w:val="1"/><w:attr w:name="Month" w:val="9"/></w:smartTagPr><w:r w:rsidR="00EC7CF2"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t>September 1</w:t></w:r><w:r w:rsidR="00D37CDD"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t>, 2008</w:t></w:r></w:smartTag><w:r w:rsidR="00406E61"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t xml:space="preserve"> &#8211; </w:t></w:r><w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date"><w:smartTagPr><w:attr w:name="Year" w:val="2008"/><w:attr w:name="Day" w:val="3"/><w:attr w:name="Month" w:val="9"/></w:smartTagPr><w:r w:rsidR="00EC7CF2"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t>September 3</w:t></w:r><w:r w:rsidR="00406E61"><w:rPr><w:i/><w:iCs/><w:sz w:val="16"/></w:rPr><w:t>, 2008</w:t></w:r></w:smartTag></w:p>
'@
$pDate.Range.InsertXML($xmlDate)

# ---------------------------------------------------------------------------
# 2) "Goal" heading (Heading 3): remove the _GoBack bookmark that sits at
#    the very start of the paragraph.
# ---------------------------------------------------------------------------
$pGoal = Get-ParaByExactText $d "Goal"
$xmlGoal = @'
<w:p w:rsidR="00EF400A" w:rsidRDefault="00EF400A" w:rsidP="00FF5F4B"><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Goa</w:t></w:r><w:r w:rsidR="00FF5F4B"><w:t>l</w:t></w:r></w:p>
'@
$pGoal.Range.InsertXML($xmlGoal)

# ---------------------------------------------------------------------------
# 3) "Time" section date range paragraph:
#    - reorder the w:attr children of both smartTagPr elements from
#      Month/Day/Year to Year/Day/Month.
#    - split the ", " run after "September 1 " into "," + a new _GoBack
#      bookmark + " ".
# ---------------------------------------------------------------------------
$pTimeRange = $null
$prevWasTimeHeading = $false
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.Trim()
    if ($prevWasTimeHeading) {
        $pTimeRange = $p
        break
    }
    $prevWasTimeHeading = ($t -eq "Time")
}
$xmlTimeRange = @'
<w:p w:rsidR="00FF5F4B" w:rsidRDefault="00EC7CF2" w:rsidP="001A67E0"><w:pPr><w:ind w:left="426"/></w:pPr><w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date"><w:smartTagPr><w:attr w:name="Year" w:val="2008"/><w:attr w:name="Day" w:val="1"/><w:attr w:name="Month" w:val="9"/></w:smartTagPr><w:r><w:t xml:space="preserve">September 1 </w:t></w:r><w:r w:rsidR="00916C2E"><w:t>,</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r w:rsidR="00916C2E"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00FF5F4B"><w:t>2008</w:t></w:r></w:smartTag><w:r w:rsidR="00916C2E"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00FF5F4B"><w:t xml:space="preserve">&#8211; </w:t></w:r><w:smartTag w:uri="urn:schemas-microsoft-com:office:smarttags" w:element="date"><w:smartTagPr><w:attr w:name="Year" w:val="2008"/><w:attr w:name="Day" w:val="3"/><w:attr w:name="Month" w:val="9"/></w:smartTagPr><w:r><w:t>September 3</w:t></w:r><w:r w:rsidR="00FE5798"><w:t>, 2008</w:t></w:r></w:smartTag></w:p>
'@
$pTimeRange.Range.InsertXML($xmlTimeRange)

# ---------------------------------------------------------------------------
# 4) Project-steps color markings: change the gray (808080) legend-style
#    markers to brown (996633) for both the "Order the ideas..." and
#    "Sub-categorize ideas" lines.
# ---------------------------------------------------------------------------
$pOrder = Get-ParaByExactText $d "/ Order the ideas by documentation order"
$xmlOrder = @'
<w:p w:rsidR="007B622F" w:rsidRPr="007829F9" w:rsidRDefault="007829F9" w:rsidP="00EC7CF2"><w:pPr><w:ind w:left="858" w:hanging="148"/><w:rPr><w:color w:val="996633"/></w:rPr></w:pPr><w:r w:rsidRPr="007829F9"><w:rPr><w:color w:val="996633"/></w:rPr><w:t>/</w:t></w:r><w:r w:rsidR="007B622F" w:rsidRPr="007829F9"><w:rPr><w:color w:val="996633"/></w:rPr><w:t xml:space="preserve"> Order the ideas by documentation order </w:t></w:r></w:p>
'@
$pOrder.Range.InsertXML($xmlOrder)

$pSubCat = Get-ParaByExactText $d "/ Sub-categorize ideas"
$xmlSubCat = @'
<w:p w:rsidR="007B622F" w:rsidRPr="007829F9" w:rsidRDefault="007829F9" w:rsidP="00EC7CF2"><w:pPr><w:ind w:left="858" w:hanging="148"/><w:rPr><w:color w:val="996633"/></w:rPr></w:pPr><w:r w:rsidRPr="007829F9"><w:rPr><w:color w:val="996633"/></w:rPr><w:t>/</w:t></w:r><w:r w:rsidR="007B622F" w:rsidRPr="007829F9"><w:rPr><w:color w:val="996633"/></w:rPr><w:t xml:space="preserve"> Sub-categorize ideas</w:t></w:r></w:p>
'@
$pSubCat.Range.InsertXML($xmlSubCat)

Write-Host "Done."
